$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(1)
$p.Range.ParagraphFormat.Hyphenation = $true
Write-Output "done"
